# Apply the "Add Purchase service, fix inventory update" data edits to
# Sheet1 of the products workbook:
#   - remove the stray "Test Product" scratch row (row 6)
#   - refresh product IDs / descriptions / prices / quantities for the
#     widget catalog rows
#   - move the active selection to E5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Widget B): price is now stored as text "20" ------------------
$ws.Range("D3").Value = "'20"
$ws.Range("D3").Style = "Normal"

# --- Row 4 (Widget C): refreshed product id -------------------------------
$ws.Range("A4").Value = "23dff0ee-e4b4-4fdf-a9ed-52c93a98eaef"

# --- Row 5 (Widget D): refreshed product id -------------------------------
$ws.Range("A5").Value = "cce54d6e-76c7-46ea-96ba-30b6be39eb8f"

# --- Row 5 (Widget D): description casing fix -----------------------------
$ws.Range("C5").Value = "A good-qauality widget"

# --- Row 3 (Widget B): description casing fix -----------------------------
$ws.Range("C3").Value = "A dazzling widget"

# --- Row 4 (Widget C): description casing fix -----------------------------
$ws.Range("C4").Value = "A premium widget"

# --- Remaining numeric corrections ----------------------------------------
$ws.Range("E2").Value = 45
$ws.Range("E3").Value = 80
$ws.Range("D4").Value = 35
$ws.Range("E4").Value = 30
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 25

# --- Drop the leftover "Test Product" scratch row (row 6) ----------------
$ws.Rows.Item(6).Delete()

# --- Match the author's final selection state -----------------------------
$ws.Range("E5").Select()
